$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Foglio1")
$ws1.Range("O9").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 7
$win.ScrollRow = 1
